$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-19 -> 2023-09-20) for every data row (rows 2-135).
$ws.Range("C2:C135").Value = 45189
